$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values are plain text (including percent strings and multi-dot
# "numbers" that Excel cannot parse as numeric literals) -- direct assignment is safe.
$ws.Range("D2").Value = "27.522.92"
$ws.Range("E2").Value = "  +5.48%  "
$ws.Range("D3").Value = "1.724.87"
$ws.Range("E3").Value = "  +4.47%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +3.48%  "
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +6.66%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.738.88"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("D14").Value = "1.961.81"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("E15").Value = "  +4.44%  "
$ws.Range("D16").Value = "0.0₅8300"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  +4.10%  "
$ws.Range("D18").Value = "27.535.18"
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("E19").Value = "  +15.23%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("E26").Value = "  +14.81%  "
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  +4.67%  "
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  +6.52%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.056.83"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "1.868.91"
$ws.Range("E45").Value = "  +4.33%  "
$ws.Range("E46").Value = "  +6.76%  "
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E51").Value = "  +1.69%  "

# --- Cells whose new values are valid numeric literals (e.g. "1.004") but must remain
# stored as Text (matching the source data which keeps these as inline strings).
# Temporarily mark the cell as Text-formatted so Excel does not coerce the literal into
# a number, then restore the default "Normal" style so no visible formatting change is
# left on the cell (NumberFormat reverts along with the style).
$textCells = @("D4", "D5", "D7", "D9", "D10", "D11", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D4").Value = "1.004"
$ws.Range("D5").Value = "225.77"
$ws.Range("D7").Value = "1.004"
$ws.Range("D9").Value = "0.06601"
$ws.Range("D10").Value = "21.72"
$ws.Range("D11").Value = "0.07709"
$ws.Range("D15").Value = "0.5843"
$ws.Range("D17").Value = "68.01"
$ws.Range("D19").Value = "219.89"
$ws.Range("D20").Value = "1.004"
$ws.Range("D21").Value = "4.723"
$ws.Range("D22").Value = "10.66"
$ws.Range("D23").Value = "6.096"
$ws.Range("D25").Value = "148.39"
$ws.Range("D26").Value = "1.726"
$ws.Range("D27").Value = "0.1233"
$ws.Range("D28").Value = "7.409"
$ws.Range("D29").Value = "16.64"
$ws.Range("D30").Value = "0.05564"
$ws.Range("D33").Value = "3.458"
$ws.Range("D35").Value = "2.841"
$ws.Range("D36").Value = "0.9599"
$ws.Range("D37").Value = "2.431"
$ws.Range("D38").Value = "0.5951"
$ws.Range("D39").Value = "0.01651"
$ws.Range("D40").Value = "5.931"
$ws.Range("D41").Value = "0.8551"
$ws.Range("D44").Value = "101.62"
$ws.Range("D47").Value = "59.00"
$ws.Range("D48").Value = "0.4437"
$ws.Range("D49").Value = "8.168"
$ws.Range("D50").Value = "1.001"
$ws.Range("D51").Value = "0.05253"
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
